{"js": "// The package-description paragraphs were reorganised: the big paragraph that\n// mixed the \"user\" / \"content\" / \"usefull_interfaces\" blurbs (joined with\n// manual line breaks) got split into separate paragraphs, and the paragraphs\n// that already existed for \"usefull_interfaces regroupe...\", \"Dans\n// usefull_interfaces...\" and \"Le package genres...\" were shifted down / had\n// their text swapped, plus one new sentence was appended to the\n// \"usefull_interfaces, langues, et genres\" paragraph.\n\nconst TEXT_A =\n  \"Le package user contient des classes li\u00e9es directement \u00e0 l\\u2019utilisateur, tel la classe Utilisateur, et la classe Marathon, li\u00e9 intimement \u00e0 un Utilisateur.\";\nconst TEXT_B =\n  \"Le package content rassemble tout le contenu de remplissage de la plateforme, donc tout les contenus Vid\u00e9oludiques (films, s\u00e9ries animes\\u2026). Le sous-package episodique repr\u00e9sente tous les contenus vid\u00e9oludiques sous forme d\\u2019\u00e9pisodes (s\u00e9ries et animes), ainsi que leur composants (\u00e9pisodes et saisons).\";\nconst TEXT_C =\n  \"Le package usefull_interfaces, langues, et genres sont tous trois des packages de services, ils servent aux classes des packages user et content. Ils \u00e0 la racine du projet de mod\u00e9lisation, dans un souci logique. En effet ce sont des ressources que nous pourrions utiliser dans d\\u2019autre contextes que content ou user, expliquant ces namespaces.\";\nconst TEXT_D =\n  \"usefull_interfaces regroupe toutes les interfaces pouvant \u00eatre utilis\u00e9es par toutes les classes, genres regroupe deux enums pour repr\u00e9senter les genres des contenus vid\u00e9oludiques (action, shonen\\u2026) et langues regroupe les langues pour les contenus vid\u00e9oludiques.\";\nconst TEXT_E =\n  \"Dans usefull_interfaces, une interface IEstDescriptible est pr\u00e9sente, et est import\u00e9 par ContenuVideoludique et Episode qui l\\u2019impl\u00e9mentent, car sont descriptibles. Ces deux classes impl\u00e9mentent aussi IEstAjoutableAuMarathon, pour pouvoir les ajouter \u00e0 la liste de lecture d\\u2019un marathon. usefull-interfaces permet ainsi de regrouper toutes interfaces susceptibles d\\u2019\u00eatre utilis\u00e9e dans tout les packages.\";\nconst TEXT_F =\n  \"Le package genres est import\u00e9 par content, pour d\u00e9crire les genres Globaux de ContenuVideoludique, mais est aussi import\u00e9 par episodique, pour d\u00e9crire les genres d\\u2019animes de la classe Anime.\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the four paragraphs we need to rewrite by their (unique) current\n// starting text, instead of trusting fixed indices.\nconst findStartingWith = (prefix) => {\n  const hit = paragraphs.items.find((p) => p.text.trim().startsWith(prefix));\n  if (!hit) {\n    throw new Error(\"Paragraph starting with '\" + prefix + \"' not found\");\n  }\n  return hit;\n};\n\nconst pUser = findStartingWith(\"Le package user contient\");\nconst pInterfacesRegroupe = findStartingWith(\"usefull_interfaces regroupe\");\nconst pDansInterfaces = findStartingWith(\"Dans usefull_interfaces\");\nconst pGenresImporte = findStartingWith(\"Le package genres est import\u00e9\");\n\n// 1) The \"user\" paragraph keeps only its first sentence.\npUser.clear();\npUser.insertText(TEXT_A, Word.InsertLocation.start);\n\n// 2) The paragraph that used to read \"usefull_interfaces regroupe...\" now\n//    holds the \"Le package content rassemble...\" text.\npInterfacesRegroupe.clear();\npInterfacesRegroupe.insertText(TEXT_B, Word.InsertLocation.start);\n\n// 3) The paragraph that used to read \"Dans usefull_interfaces...\" now holds\n//    the \"Le package usefull_interfaces, langues, et genres...\" text, with a\n//    new closing sentence appended.\npDansInterfaces.clear();\npDansInterfaces.insertText(TEXT_C, Word.InsertLocation.start);\n\n// 4) The paragraph that used to read \"Le package genres est import\u00e9...\" now\n//    holds the old \"usefull_interfaces regroupe...\" text.\npGenresImporte.clear();\npGenresImporte.insertText(TEXT_D, Word.InsertLocation.start);\n\n// 5) Two new paragraphs are inserted right after it: the old \"Dans\n//    usefull_interfaces...\" text, followed by the old \"Le package genres est\n//    import\u00e9...\" text. Both inherit the indentation of the reference\n//    paragraph automatically.\nconst pE = pGenresImporte.insertParagraph(TEXT_E, Word.InsertLocation.after);\npE.insertParagraph(TEXT_F, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# The package-description paragraphs were reorganised: the big paragraph that\n# mixed the \"user\" / \"content\" / \"usefull_interfaces\" blurbs (joined with\n# manual line breaks) got split into separate paragraphs, and the paragraphs\n# that already existed for \"usefull_interfaces regroupe...\", \"Dans\n# usefull_interfaces...\" and \"Le package genres...\" were shifted down / had\n# their text swapped, plus one new sentence was appended to the\n# \"usefull_interfaces, langues, et genres\" paragraph.\n\n$TEXT_A = \"Le package user contient des classes li\u00e9es directement \u00e0 l\u2019utilisateur, tel la classe Utilisateur, et la classe Marathon, li\u00e9 intimement \u00e0 un Utilisateur.\"\n$TEXT_B = \"Le package content rassemble tout le contenu de remplissage de la plateforme, donc tout les contenus Vid\u00e9oludiques (films, s\u00e9ries animes\u2026). Le sous-package episodique repr\u00e9sente tous les contenus vid\u00e9oludiques sous forme d\u2019\u00e9pisodes (s\u00e9ries et animes), ainsi que leur composants (\u00e9pisodes et saisons).\"\n$TEXT_C = \"Le package usefull_interfaces, langues, et genres sont tous trois des packages de services, ils servent aux classes des packages user et content. Ils \u00e0 la racine du projet de mod\u00e9lisation, dans un souci logique. En effet ce sont des ressources que nous pourrions utiliser dans d\u2019autre contextes que content ou user, expliquant ces namespaces.\"\n$TEXT_D = \"usefull_interfaces regroupe toutes les interfaces pouvant \u00eatre utilis\u00e9es par toutes les classes, genres regroupe deux enums pour repr\u00e9senter les genres des contenus vid\u00e9oludiques (action, shonen\u2026) et langues regroupe les langues pour les contenus vid\u00e9oludiques.\"\n$TEXT_E = \"Dans usefull_interfaces, une interface IEstDescriptible est pr\u00e9sente, et est import\u00e9 par ContenuVideoludique et Episode qui l\u2019impl\u00e9mentent, car sont descriptibles. Ces deux classes impl\u00e9mentent aussi IEstAjoutableAuMarathon, pour pouvoir les ajouter \u00e0 la liste de lecture d\u2019un marathon. usefull-interfaces permet ainsi de regrouper toutes interfaces susceptibles d\u2019\u00eatre utilis\u00e9e dans tout les packages.\"\n$TEXT_F = \"Le package genres est import\u00e9 par content, pour d\u00e9crire les genres Globaux de ContenuVideoludique, mais est aussi import\u00e9 par episodique, pour d\u00e9crire les genres d\u2019animes de la classe Anime.\"\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($doc, $prefix) {\n    $paras = $doc.Paragraphs\n    for ($i = 1; $i -le $paras.Count; $i++) {\n        if ($paras.Item($i).Range.Text.StartsWith($prefix)) {\n            return $i\n        }\n    }\n    return -1\n}\n\nfunction Set-ParagraphText($doc, $index, $newText) {\n    # Insert a fresh paragraph right after the target (it inherits the\n    # target's formatting/indentation), fill it with the new text, then\n    # delete the old paragraph outright. This avoids leaving behind orphan\n    # run-level markup (e.g. spell-check <w:proofErr/> tags) that a plain\n    # text-range delete on a multi-run paragraph can strand.\n    $doc.Paragraphs.Item($index).Range.InsertParagraphAfter()\n    $doc.Paragraphs.Item($index + 1).Range.InsertBefore($newText)\n    $doc.Paragraphs.Item($index).Range.Delete()\n}\n\nfunction Insert-ParagraphAfter($doc, $index, $newText) {\n    $doc.Paragraphs.Item($index).Range.InsertParagraphAfter()\n    $doc.Paragraphs.Item($index + 1).Range.InsertBefore($newText)\n    return $index + 1\n}\n\n# Locate the four paragraphs we need to rewrite by their (unique) current\n# starting text, instead of trusting fixed indices.\n$idxUser = Find-ParagraphIndex $d \"Le package user contient\"\n$idxInterfacesRegroupe = Find-ParagraphIndex $d \"usefull_interfaces regroupe\"\n$idxDansInterfaces = Find-ParagraphIndex $d \"Dans usefull_interfaces\"\n$idxGenresImporte = Find-ParagraphIndex $d \"Le package genres est import\u00e9\"\n\n# 1) The \"user\" paragraph keeps only its first sentence.\nSet-ParagraphText $d $idxUser $TEXT_A\n\n# 2) The paragraph that used to read \"usefull_interfaces regroupe...\" now\n#    holds the \"Le package content rassemble...\" text.\nSet-ParagraphText $d $idxInterfacesRegroupe $TEXT_B\n\n# 3) The paragraph that used to read \"Dans usefull_interfaces...\" now holds\n#    the \"Le package usefull_interfaces, langues, et genres...\" text, with a\n#    new closing sentence appended.\nSet-ParagraphText $d $idxDansInterfaces $TEXT_C\n\n# 4) The paragraph that used to read \"Le package genres est import\u00e9...\" now\n#    holds the old \"usefull_interfaces regroupe...\" text.\nSet-ParagraphText $d $idxGenresImporte $TEXT_D\n\n# 5) Two new paragraphs are inserted right after it: the old \"Dans\n#    usefull_interfaces...\" text, followed by the old \"Le package genres est\n#    import\u00e9...\" text. Both inherit the indentation of the reference\n#    paragraph automatically.\n$newIdx = Insert-ParagraphAfter $d $idxGenresImporte $TEXT_E\n$null = Insert-ParagraphAfter $d $newIdx $TEXT_F\n"}
